$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3084.3157  # H74: was 3084.3684
$ws.Cells.Item(74, 9).Value = 2326  # I74: was 2326.125
$ws.Cells.Item(74, 11).Value = 2326  # K74: was 2326.125
$ws.Cells.Item(74, 13).Value = -1390  # M74: was -1390.125

$ws.Cells.Item(77, 8).Value = 3084.3157  # H77: was 3084.3684
$ws.Cells.Item(77, 9).Value = 2326  # I77: was 2326.125
$ws.Cells.Item(77, 11).Value = 11630  # K77: was 11630.625
$ws.Cells.Item(77, 13).Value = -6950  # M77: was -6950.625

$ws.Cells.Item(132, 8).Value = 13501.529  # H132: was 11639.55
$ws.Cells.Item(132, 9).Value = 8885.77  # I132: was 7423.75
$ws.Cells.Item(132, 11).Value = 26657.31  # K132: was 22271.25
$ws.Cells.Item(132, 13).Value = -24127.31  # M132: was -19741.25

$ws.Cells.Item(141, 8).Value = 4818.846  # H141: was 5192
$ws.Cells.Item(141, 9).Value = 5185.5557  # I141: was 5645.7144
$ws.Cells.Item(141, 10).Value = 3993.75  # J141: was 4133.3335
$ws.Cells.Item(141, 11).Value = 15556.6671  # K141: was 16937.1432
$ws.Cells.Item(141, 12).Value = 11981.25  # L141: was 12400.0005
$ws.Cells.Item(141, 13).Value = -10376.6671  # M141: was -11757.1432
$ws.Cells.Item(141, 14).Value = -22341.25  # N141: was -22760.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5206.7935  # H32: was 5287.161
$ws.Cells.Item(32, 9).Value = 2921.75  # I32: was 2979.149
$ws.Cells.Item(32, 11).Value = 2921.75  # K32: was 2979.149
$ws.Cells.Item(32, 13).Value = -2634.75  # M32: was -2692.149

$ws.Cells.Item(45, 8).Value = 2122.7058  # H45: was 2292.4614
$ws.Cells.Item(45, 9).Value = 1757.9166  # I45: was 1978.5555
$ws.Cells.Item(45, 10).Value = 2998.2  # J45: was 2998.75
$ws.Cells.Item(45, 11).Value = 1757.9166  # K45: was 1978.5555
$ws.Cells.Item(45, 12).Value = 2998.2  # L45: was 2998.75
$ws.Cells.Item(45, 13).Value = -1380.9166  # M45: was -1601.5555
$ws.Cells.Item(45, 14).Value = -3752.2  # N45: was -3752.75

$ws.Cells.Item(133, 8).Value = 98109.55499999999  # H133: was 99248.25
$ws.Cells.Item(133, 10).Value = 98109.55499999999  # J133: was 99248.25
$ws.Cells.Item(133, 12).Value = 98109.55499999999  # L133: was 99248.25
$ws.Cells.Item(133, 14).Value = -103169.555  # N133: was -104308.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 11304.923  # H134: was 11107.275
$ws.Cells.Item(134, 9).Value = 6860.161  # I134: was 6661.4062
$ws.Cells.Item(134, 10).Value = 28528.375  # J134: was 28890.75
$ws.Cells.Item(134, 11).Value = 20580.483  # K134: was 19984.2186
$ws.Cells.Item(134, 12).Value = 85585.125  # L134: was 86672.25
$ws.Cells.Item(134, 13).Value = -18045.483  # M134: was -17449.2186
$ws.Cells.Item(134, 14).Value = -90655.125  # N134: was -91742.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 581.75  # H7: was 577.2
$ws.Cells.Item(7, 9).Value = 178.33333  # I7: was 165.15384
$ws.Cells.Item(7, 10).Value = 1186.875  # J7: was 1342.4286
$ws.Cells.Item(7, 11).Value = 178.33333  # K7: was 165.15384
$ws.Cells.Item(7, 12).Value = 1186.875  # L7: was 1342.4286
$ws.Cells.Item(7, 13).Value = -65.33332999999999  # M7: was -52.15384
$ws.Cells.Item(7, 14).Value = -1412.875  # N7: was -1568.4286

$ws.Cells.Item(11, 8).Value = 1170.1428  # H11: was 1593.25
$ws.Cells.Item(11, 9).Value = 699.5  # I11: was 999
$ws.Cells.Item(11, 10).Value = 1358.4  # J11: was 1791.3334
$ws.Cells.Item(11, 11).Value = 699.5  # K11: was 999
$ws.Cells.Item(11, 12).Value = 1358.4  # L11: was 1791.3334
$ws.Cells.Item(11, 13).Value = -559.5  # M11: was -859
$ws.Cells.Item(11, 14).Value = -1638.4  # N11: was -2071.3334

$ws.Cells.Item(31, 8).Value = 66505.766  # H31: was 55104.523
$ws.Cells.Item(31, 9).Value = 204634.7  # I31: was 114107.555
$ws.Cells.Item(31, 11).Value = 204634.7  # K31: was 114107.555
$ws.Cells.Item(31, 13).Value = -204339.7  # M31: was -113812.555

$ws.Cells.Item(34, 8).Value = 66505.766  # H34: was 55104.523
$ws.Cells.Item(34, 9).Value = 204634.7  # I34: was 114107.555
$ws.Cells.Item(34, 11).Value = 204634.7  # K34: was 114107.555
$ws.Cells.Item(34, 13).Value = -204432.7  # M34: was -113905.555

$ws.Cells.Item(55, 8).Value = 0  # H55: was 9000
$ws.Cells.Item(55, 9).Value = 0  # I55: was 9000
$ws.Cells.Item(55, 11).Value = 0  # K55: was 9000
$ws.Cells.Item(55, 13).ClearContents()  # M55: was -8685

$ws.Cells.Item(58, 8).Value = 14061.637  # H58: was 12634
$ws.Cells.Item(58, 9).Value = 7131.5835  # I58: was 5562.6875
$ws.Cells.Item(58, 11).Value = 7131.5835  # K58: was 5562.6875
$ws.Cells.Item(58, 13).Value = -6928.5835  # M58: was -5359.6875

$ws.Cells.Item(86, 8).Value = 7838.591  # H86: was 8069.048
$ws.Cells.Item(86, 9).Value = 8066  # I86: was 8403.799999999999
$ws.Cells.Item(86, 11).Value = 8066  # K86: was 8403.799999999999
$ws.Cells.Item(86, 13).Value = -6943  # M86: was -7280.799999999999

$ws.Cells.Item(89, 8).Value = 7838.591  # H89: was 8069.048
$ws.Cells.Item(89, 9).Value = 8066  # I89: was 8403.799999999999
$ws.Cells.Item(89, 11).Value = 40330  # K89: was 42019
$ws.Cells.Item(89, 13).Value = -34714  # M89: was -36403

$ws.Cells.Item(99, 8).Value = 6463.967  # H99: was 6564.2666
$ws.Cells.Item(99, 9).Value = 4443.4375  # I99: was 4631.5
$ws.Cells.Item(99, 11).Value = 4443.4375  # K99: was 4631.5
$ws.Cells.Item(99, 13).Value = -2945.4375  # M99: was -3133.5

$ws.Cells.Item(107, 8).Value = 4060.9412  # H107: was 4072.7058
$ws.Cells.Item(107, 10).Value = 8707  # J107: was 8740.333000000001
$ws.Cells.Item(107, 12).Value = 8707  # L107: was 8740.333000000001
$ws.Cells.Item(107, 14).Value = -12547  # N107: was -12580.333

$ws.Cells.Item(122, 8).Value = 7873.222  # H122: was 7235.8
$ws.Cells.Item(122, 9).Value = 2183.1667  # I122: was 2085.4285
$ws.Cells.Item(122, 11).Value = 6549.500100000001  # K122: was 6256.2855
$ws.Cells.Item(122, 13).Value = -4099.500100000001  # M122: was -3806.2855

$ws.Cells.Item(126, 8).Value = 6463.967  # H126: was 6564.2666
$ws.Cells.Item(126, 9).Value = 4443.4375  # I126: was 4631.5
$ws.Cells.Item(126, 11).Value = 13330.3125  # K126: was 13894.5
$ws.Cells.Item(126, 13).Value = -10860.3125  # M126: was -11424.5

$ws.Cells.Item(133, 8).Value = 121888.555  # H133: was 124999.06
$ws.Cells.Item(133, 10).Value = 121888.555  # J133: was 124999.06
$ws.Cells.Item(133, 12).Value = 121888.555  # L133: was 124999.06
$ws.Cells.Item(133, 14).Value = -126948.555  # N133: was -130059.06

$ws.Cells.Item(134, 8).Value = 41675484  # H134: was 55566220
$ws.Cells.Item(134, 9).Value = 3640.8333  # I134: was 4005
$ws.Cells.Item(134, 11).Value = 10922.4999  # K134: was 12015
$ws.Cells.Item(134, 13).Value = -8387.499899999999  # M134: was -9480

$ws.Cells.Item(136, 8).Value = 14061.637  # H136: was 12634
$ws.Cells.Item(136, 9).Value = 7131.5835  # I136: was 5562.6875
$ws.Cells.Item(136, 11).Value = 21394.7505  # K136: was 16688.0625
$ws.Cells.Item(136, 13).Value = -18844.7505  # M136: was -14138.0625

$ws.Cells.Item(137, 8).Value = 54750  # H137: was 54888.332
$ws.Cells.Item(137, 10).Value = 54750  # J137: was 54888.332
$ws.Cells.Item(137, 12).Value = 54750  # L137: was 54888.332
$ws.Cells.Item(137, 14).Value = -64950  # N137: was -65088.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 312.375  # H18: was 266.16666
$ws.Cells.Item(18, 9).Value = 312.375  # I18: was 266.16666
$ws.Cells.Item(18, 11).Value = 937.125  # K18: was 798.4999799999999
$ws.Cells.Item(18, 13).Value = -768.125  # M18: was -629.4999799999999

$ws.Cells.Item(112, 8).Value = 6807.25  # H112: was 10608.0625
$ws.Cells.Item(112, 10).Value = 10965  # J112: was 11745
$ws.Cells.Item(112, 12).Value = 32895  # L112: was 35235
$ws.Cells.Item(112, 14).Value = -35111  # N112: was -37451

$ws.Cells.Item(129, 8).Value = 2343  # H129: was 2463.2222
$ws.Cells.Item(129, 9).Value = 784.125  # I129: was 833.2857
$ws.Cells.Item(129, 10).Value = 2936.8572  # J129: was 3033.7
$ws.Cells.Item(129, 11).Value = 2352.375  # K129: was 2499.8571
$ws.Cells.Item(129, 12).Value = 8810.571599999999  # L129: was 9101.099999999999
$ws.Cells.Item(129, 13).Value = 2647.625  # M129: was 2500.1429
$ws.Cells.Item(129, 14).Value = -18810.5716  # N129: was -19101.1

$ws.Cells.Item(131, 8).Value = 1433.83  # H131: was 1714.29
$ws.Cells.Item(131, 9).Value = 780.25  # I131: was 860.2222
$ws.Cells.Item(131, 10).Value = 1490.6631  # J131: was 1798.7583
$ws.Cells.Item(131, 11).Value = 2340.75  # K131: was 2580.6666
$ws.Cells.Item(131, 12).Value = 4471.9893  # L131: was 5396.2749
$ws.Cells.Item(131, 13).Value = 2699.25  # M131: was 2459.3334
$ws.Cells.Item(131, 14).Value = -14551.9893  # N131: was -15476.2749

$ws.Cells.Item(132, 8).Value = 2756475.5  # H132: was 10102073
$ws.Cells.Item(132, 9).Value = 1642.8572  # I132: was 0
$ws.Cells.Item(132, 10).Value = 7577432.5  # J132: was 10102073
$ws.Cells.Item(132, 11).Value = 14785.7148  # K132: was 0
$ws.Cells.Item(132, 12).Value = 68196892.5  # L132: was 90918657
$ws.Cells.Item(132, 13).Value = -12255.7148  # M132: newly added (was empty)
$ws.Cells.Item(132, 14).Value = -68201952.5  # N132: was -90923717

$ws.Cells.Item(138, 8).Value = 2802.2856  # H138: was 3207.7778
$ws.Cells.Item(138, 9).Value = 2604.3333  # I138: was 2604.6667
$ws.Cells.Item(138, 10).Value = 3990  # J138: was 3380.0952
$ws.Cells.Item(138, 11).Value = 7812.999899999999  # K138: was 7814.000100000001
$ws.Cells.Item(138, 12).Value = 11970  # L138: was 10140.2856
$ws.Cells.Item(138, 13).Value = -2672.999899999999  # M138: was -2674.000100000001
$ws.Cells.Item(138, 14).Value = -22250  # N138: was -20420.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 12482.131  # H7: was 12395.174
$ws.Cells.Item(7, 9).Value = 13644.818  # I7: was 13553.909
$ws.Cells.Item(7, 10).Value = 11416.333  # J7: was 11333
$ws.Cells.Item(7, 11).Value = 13644.818  # K7: was 13553.909
$ws.Cells.Item(7, 12).Value = 11416.333  # L7: was 11333
$ws.Cells.Item(7, 13).Value = -13532.818  # M7: was -13441.909
$ws.Cells.Item(7, 14).Value = -11640.333  # N7: was -11557

$ws.Cells.Item(22, 8).Value = 11891.8  # H22: was 9236.857
$ws.Cells.Item(22, 9).Value = 2080  # I22: was 2339.75
$ws.Cells.Item(22, 11).Value = 2080  # K22: was 2339.75
$ws.Cells.Item(22, 13).Value = -1785  # M22: was -2044.75

$ws.Cells.Item(27, 8).Value = 11891.8  # H27: was 9236.857
$ws.Cells.Item(27, 9).Value = 2080  # I27: was 2339.75
$ws.Cells.Item(27, 11).Value = 2080  # K27: was 2339.75
$ws.Cells.Item(27, 13).Value = -1973  # M27: was -2232.75

$ws.Cells.Item(40, 8).Value = 7966.923  # H40: was 7925.76
$ws.Cells.Item(40, 9).Value = 5619.077  # I40: was 5289.143
$ws.Cells.Item(40, 10).Value = 10314.77  # J40: was 11281.454
$ws.Cells.Item(40, 11).Value = 5619.077  # K40: was 5289.143
$ws.Cells.Item(40, 12).Value = 10314.77  # L40: was 11281.454
$ws.Cells.Item(40, 13).Value = -5483.077  # M40: was -5153.143
$ws.Cells.Item(40, 14).Value = -10586.77  # N40: was -11553.454

$ws.Cells.Item(61, 8).Value = 4324.9546  # H61: was 3170.121
$ws.Cells.Item(61, 9).Value = 3964.2  # I61: was 2729.12
$ws.Cells.Item(61, 10).Value = 5098  # J61: was 4548.25
$ws.Cells.Item(61, 11).Value = 3964.2  # K61: was 2729.12
$ws.Cells.Item(61, 12).Value = 5098  # L61: was 4548.25
$ws.Cells.Item(61, 13).Value = -3762.2  # M61: was -2527.12
$ws.Cells.Item(61, 14).Value = -5502  # N61: was -4952.25

$ws.Cells.Item(113, 8).Value = 4324.9546  # H113: was 3170.121
$ws.Cells.Item(113, 9).Value = 3964.2  # I113: was 2729.12
$ws.Cells.Item(113, 10).Value = 5098  # J113: was 4548.25
$ws.Cells.Item(113, 11).Value = 3964.2  # K113: was 2729.12
$ws.Cells.Item(113, 12).Value = 5098  # L113: was 4548.25
$ws.Cells.Item(113, 13).Value = -1794.2  # M113: was -559.1199999999999
$ws.Cells.Item(113, 14).Value = -9438  # N113: was -8888.25

$ws.Cells.Item(126, 8).Value = 12482.131  # H126: was 12395.174
$ws.Cells.Item(126, 9).Value = 13644.818  # I126: was 13553.909
$ws.Cells.Item(126, 10).Value = 11416.333  # J126: was 11333
$ws.Cells.Item(126, 11).Value = 40934.454  # K126: was 40661.727
$ws.Cells.Item(126, 12).Value = 34248.999  # L126: was 33999
$ws.Cells.Item(126, 13).Value = -38464.454  # M126: was -38191.727
$ws.Cells.Item(126, 14).Value = -39188.999  # N126: was -38939

$ws.Cells.Item(136, 8).Value = 732031.7  # H136: was 825464.3
$ws.Cells.Item(136, 9).Value = 19005.46  # I136: was 23829.3
$ws.Cells.Item(136, 10).Value = 1031042.6  # J136: was 1101890.1
$ws.Cells.Item(136, 11).Value = 57016.38  # K136: was 71487.89999999999
$ws.Cells.Item(136, 12).Value = 3093127.8  # L136: was 3305670.3
$ws.Cells.Item(136, 13).Value = -54466.38  # M136: was -68937.89999999999
$ws.Cells.Item(136, 14).Value = -3098227.8  # N136: was -3310770.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 6835.6665  # H113: was 6840.1113
$ws.Cells.Item(113, 9).Value = 9513.833000000001  # I113: was 11136.6
$ws.Cells.Item(113, 10).Value = 1479.3334  # J113: was 1469.5
$ws.Cells.Item(113, 11).Value = 28541.499  # K113: was 33409.8
$ws.Cells.Item(113, 12).Value = 4438.0002  # L113: was 4408.5
$ws.Cells.Item(113, 13).Value = -26371.499  # M113: was -31239.8
$ws.Cells.Item(113, 14).Value = -8778.0002  # N113: was -8748.5
